$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells with the same style as the existing header row (copy format from E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# Fill boolean FALSE values for rows 2 through 20 in columns F, G, H
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
